$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.747.25"
$ws.Range("E2").Value = "  +7.04%  "
$ws.Range("D3").Value = "2.618.23"
$ws.Range("E3").Value = "  +7.12%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "186.86"
$ws.Range("E5").Value = "  +13.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "581.60"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  +4.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.197"
$ws.Range("E9").Value = "  +14.94%  "
$ws.Range("D10").Value = "2.617.90"
$ws.Range("E10").Value = "  +7.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.163"
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +7.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.68"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "73.855.78"
$ws.Range("E14").Value = "  +7.35%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").Value = "  +4.87%  "
$ws.Range("D16").Value = "3.103.06"
$ws.Range("E16").Value = "  +7.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.50"
$ws.Range("E17").Value = "  +12.79%  "
$ws.Range("D18").Value = "2.619.11"
$ws.Range("E18").Value = "  +7.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.17"
$ws.Range("E19").Value = "  +30.45%  "
$ws.Range("E20").Value = "  +10.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.82"
$ws.Range("E21").Value = "  +7.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.28"
$ws.Range("E22").Value = "  +17.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.07"
$ws.Range("E23").Value = "  +5.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.997"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.80"
$ws.Range("E25").Value = "  +6.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.11"
$ws.Range("E26").Value = "  +7.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("E27").Value = "  +10.35%  "
$ws.Range("D28").Value = "2.753.52"
$ws.Range("E28").Value = "  +7.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +5.12%  "
$ws.Range("D30").Value = "0.0₃0936"
$ws.Range("E30").Value = "  +12.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "523.17"
$ws.Range("E31").Value = "  +20.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.38"
$ws.Range("E32").Value = "  +13.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.63"
$ws.Range("E33").Value = "  +5.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.74"
$ws.Range("E34").Value = "  +8.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.33"
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("E37").Value = "  +9.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.08"
$ws.Range("E38").Value = "  +6.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.23"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.89"
$ws.Range("E41").Value = "  +11.38%  "
$ws.Range("E42").Value = "  +9.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.324"
$ws.Range("E43").Value = "  +7.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "162.65"
$ws.Range("E44").Value = "  +25.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.36"
$ws.Range("E45").Value = "  +12.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  +8.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.94"
$ws.Range("E47").Value = "  +3.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0852"
$ws.Range("E48").Value = "  +18.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.59"
$ws.Range("E49").Value = "  +7.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.521"
$ws.Range("E50").Value = "  +7.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.66"
$ws.Range("E51").Value = "  +21.25%  "
